# Rename several pf_* result column headers in the shared strings /
# header row of every worksheet in the workbook.
#
# Sheets 1-8 ("LLL_*", "LL_*" short variants without phase columns) use a
# shorter header row (columns A:Q).
# Sheets 9-32 use a longer header row with per-phase columns (A:AQ).
#
# The mapping below reflects the renames performed in the diff:
#   *_deg            -> *_degree
#   pf_q_*_mw        -> pf_q_*_mvar   (phase reactive power columns)
#   pf_vm_*_bus_pu   -> pf_vm_*_pu    (drop redundant "_bus")
#   pf_va_*_bus_deg  -> pf_va_*_degree

$wb = $excel.ActiveWorkbook

# Map: old header text -> new header text
$renames = @{
    "pf_ikss_from_deg"      = "pf_ikss_from_degree"
    "pf_ikss_to_deg"        = "pf_ikss_to_degree"
    "pf_va_from_deg"        = "pf_va_from_degree"
    "pf_va_to_deg"          = "pf_va_to_degree"
    "pf_q_a_from_mw"        = "pf_q_a_from_mvar"
    "pf_q_b_from_mw"        = "pf_q_b_from_mvar"
    "pf_q_c_from_mw"        = "pf_q_c_from_mvar"
    "pf_q_a_to_mw"          = "pf_q_a_to_mvar"
    "pf_q_b_to_mw"          = "pf_q_b_to_mvar"
    "pf_q_c_to_mw"          = "pf_q_c_to_mvar"
    "pf_ikss_a_from_deg"    = "pf_ikss_a_from_degree"
    "pf_ikss_b_from_deg"    = "pf_ikss_b_from_degree"
    "pf_ikss_c_from_deg"    = "pf_ikss_c_from_degree"
    "pf_ikss_a_to_deg"      = "pf_ikss_a_to_degree"
    "pf_ikss_b_to_deg"      = "pf_ikss_b_to_degree"
    "pf_ikss_c_to_deg"      = "pf_ikss_c_to_degree"
    "pf_vm_b_from_bus_pu"   = "pf_vm_b_from_pu"
    "pf_vm_c_from_bus_pu"   = "pf_vm_c_from_pu"
    "pf_vm_a_to_bus_pu"     = "pf_vm_a_to_pu"
    "pf_vm_b_to_bus_pu"     = "pf_vm_b_to_pu"
    "pf_vm_c_to_bus_pu"     = "pf_vm_c_to_pu"
    "pf_va_a_from_bus_deg"  = "pf_va_a_from_degree"
    "pf_va_b_from_bus_deg"  = "pf_va_b_from_degree"
    "pf_va_c_from_bus_deg"  = "pf_va_c_from_degree"
    "pf_va_a_to_bus_deg"    = "pf_va_a_to_degree"
    "pf_va_b_to_bus_deg"    = "pf_va_b_to_degree"
    "pf_va_c_to_bus_deg"    = "pf_va_c_to_degree"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $headerRow = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $used.Columns.Count))
    foreach ($cell in $headerRow.Cells) {
        $val = $cell.Value2
        if ($null -ne $val -and $renames.ContainsKey([string]$val)) {
            $cell.Value2 = $renames[[string]$val]
        }
    }
}
